$d = $word.ActiveDocument

# 1. "by " / "Agresti" / ", Franklin and Klingenberg, 5" were three separate runs
#    (split apart by spell-check proofErr markers around "Agresti"). Re-saving
#    collapses them into a single run with the combined text and drops the
#    now-unnecessary w:proofErr spell-check bookmarks.
[void]$d.Content.Find.Execute(
    "by Agresti, Franklin and Klingenberg, 5", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "by Agresti, Franklin and Klingenberg, 5", 2)

# 2. Fix the title mistake: "Example 13" should have been "Example 14".
[void]$d.Content.Find.Execute(
    "Example 13:", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Example 14:", 2)
